# The resume's "Software Development Engineer / Intern at Amazon" bullet
# used to be split across two runs: a plain run for
# "Software Development Engineer –" and an underlined run for
# "Intern at Amazon". The author merged them into a single underlined
# run reading "Software Development Engineer –Intern at Amazon", and -
# since this was the spot last edited - Word's automatic "_GoBack"
# bookmark moved here (it used to sit around the "Virtual "/"Machines,
# which are light," boundary further down the document).

$d = $word.ActiveDocument
$dash = [char]0x2013
$target = "Software Development Engineer " + $dash + "Intern at Amazon"

# 1. Merge the two runs into one run of text via Find & Replace (Word
#    collapses the matched range into a single run when it is replaced).
$merge = $d.Content
$merge.Find.Execute($target, $true, $false, $false, $false, $false, $true, `
                     1, $false, $target, 2) | Out-Null

# 2. Re-find the (now single-run) text and give it the underline
#    formatting that the second of the original two runs had.
$merged = $d.Content
$merged.Find.Execute($target, $true, $false, $false, $false, $false, $true, `
                      1, $false, "", 0) | Out-Null
$merged.Font.Underline = 1

# 3. Word keeps only one "_GoBack" bookmark in the whole document,
#    marking the location of the last edit. Adding it here moves it
#    off its old spot (around "Virtual "/"Machines, which are light,")
#    and wraps it around the text we just edited.
$d.Bookmarks.Add("_GoBack", $merged) | Out-Null

Write-Output "Done"
